$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header row labels
$ws.Range("A1").Value = "wafer name"
$ws.Range("B1").Value = "peak position"

# Update the wafer label for the first group of rows (row 2 only)
$ws.Range("A2").Value = "P01 B"

# Update the active cell selection to match the authored state
$ws.Range("A3").Select()
